$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($countryName, $b, $c, $d, $e, $f, $g, $h) {
    $found = $ws.Range("A4:A219").Find($countryName)
    $r = $found.Row()
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

Set-CountryRow "Estados Unidos" 6747314 36244 4021492 2526885 0 417 198937
Set-CountryRow "Brasil" 4349544 19089 3613184 604243 0 454 132117
Set-CountryRow "Colombia" 721892 5573 606925 91844 0 199 23123
Set-CountryRow "Alemania" 263221 1923 237550 16235 0 8 9436
Set-CountryRow "Canada" 138010 1351 121224 7607 0 8 9179
Set-CountryRow "Panama" 102204 459 74782 25249 0 7 2173
Set-CountryRow "Egipto" 101177 168 84969 10547 0 13 5661
Set-CountryRow "Guatemala" 82172 263 71352 7848 0 15 2972
Set-CountryRow "Japon" 75657 439 67242 6973 0 3 1442
Set-CountryRow "Venezuela" 60540 0 48644 11411 0 0 485
Set-CountryRow "Barein" 60965 658 54204 6548 0 1 213
Set-CountryRow "Nigeria" 56388 132 44337 10968 0 1 1083
Set-CountryRow "Chequia" 37222 1034 22020 14737 0 9 465
Set-CountryRow "Kenia" 36205 48 23243 12338 0 2 624
Set-CountryRow "Camerun" 20228 61 18837 976 0 0 415
Set-CountryRow "Noruega" 12276 122 10371 1640 0 0 265
Set-CountryRow "Consejo Danes para los Refugiados" 10390 0 9807 319 0 0 264
Set-CountryRow "Zimbabue" 7531 5 5690 1617 0 0 224
Set-CountryRow "Mauritania" 7295 19 6835 299 0 0 161
Set-CountryRow "Luxemburgo" 7238 0 6397 717 0 0 124
Set-CountryRow "Montenegro" 6900 188 4589 2190 0 1 121
Set-CountryRow "Tunez" 7382 747 2175 5090 0 10 117
Set-CountryRow "Malaui" 5697 7 3742 1777 0 1 178
Set-CountryRow "Uruguay" 1812 4 1528 239 0 0 45
Set-CountryRow "Comoras" 457 1 427 23 0 0 7

# Re-sort the data range by Casos totales (column B) descending, matching the source pipeline
$rng = $ws.Range("A4:H219")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B4:B219"), 0, 2)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 01:22"
